$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.55739999999998
$ws.Range("A7").Value = -21.99750000000001
$ws.Range("B7").Value = 4.632000000000002
$ws.Range("B15").Value = 4.911299999999998
$ws.Range("A16").Value = -21.55749999999998
$ws.Range("C16").Value = -11.98319999999999
$ws.Range("C19").Value = -12.27340000000001
$ws.Range("B21").Value = 10.2121
$ws.Range("B22").Value = 10.13600000000001
$ws.Range("B23").Value = 8.960100000000006
$ws.Range("A28").Value = -22.0834
$ws.Range("A29").Value = -21.30439999999998
$ws.Range("A32").Value = -21.18409999999999
$ws.Range("B34").Value = 9.636600000000005
$ws.Range("C36").Value = -12.8089
$ws.Range("A40").Value = -19.95829999999999
$ws.Range("B43").Value = 6.058000000000001
$ws.Range("B45").Value = 5.116000000000002
$ws.Range("C46").Value = -14.81479999999998
$ws.Range("B50").Value = 4.966299999999997
$ws.Range("C50").Value = -13.13409999999999
$ws.Range("B51").Value = 5.955400000000001
$ws.Range("A52").Value = -22.07409999999999
$ws.Range("A57").Value = -22.37820000000001
$ws.Range("A66").Value = -21.5156
$ws.Range("B66").Value = 5.514799999999999
$ws.Range("B67").Value = 5.190899999999999
$ws.Range("B79").Value = 9.722200000000004
$ws.Range("B84").Value = 5.557099999999999
$ws.Range("B92").Value = 4.681999999999999
$ws.Range("C95").Value = -12.1606
$ws.Range("B97").Value = 6.548399999999995
$ws.Range("C97").Value = -12.668
$ws.Range("A100").Value = -22.0006
